# Updates loading_percent values for the "380 kV" case (Case_3_119)
# on Sheet1, rows 2-25, columns B, C, D, F, G, H, I, K.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cols = @("B","C","D","F","G","H","I","K")

$data = @(
    @(9.328343784355354,5.11983886738552,5.174628543835425,24.12854153961361,28.44668901693521,14.32345971147685,20.44073493203962,9.15117348068347),
    @(8.973395608549719,4.862255450572296,5.109305836660659,24.2332812560475,28.62934061186657,14.39305911169561,20.57170718249737,8.903834857517133),
    @(8.749106100296851,4.695928554915776,5.068402551887662,24.30549255657003,28.7533022866289,14.43859279322465,20.65690240593316,8.749533716590612),
    @(8.656245283022423,4.626136008100352,5.051544776413753,24.33689655656531,28.8067674019039,14.45785202945058,20.69282180025974,8.686134418262089),
    @(8.640741828275512,4.614426797007385,5.048734430933693,24.34223031085235,28.81582279051327,14.46109252390221,20.6988587737474,8.675578372782681),
    @(8.747859472317725,4.694995398211681,5.068175954666987,24.30590808850319,28.75401142008677,14.43884968002047,20.65738196164335,8.74868067309286),
    @(9.207352480872732,5.032742108883034,5.152276490195614,24.16301046117747,28.50720313903543,14.34687655245598,20.48490268500834,9.066449723838261),
    @(10.05279978417941,5.628797866551795,5.310441866279359,23.94588173629342,28.11785732994641,14.18873070768012,20.18456932544136,9.666530447971288),
    @(10.63399157203278,6.024828065842503,5.421948401747398,23.82535048872824,27.89073457329084,14.08608476038009,19.98700301046187,10.08874630235118),
    @(10.88868502354561,6.195680978612238,5.471539063337493,23.77908969825908,27.8004633485101,14.04233053094717,19.90213703464215,10.27591875887178),
    @(10.98366586369991,6.25902770233739,5.490145312654406,23.762812136384,27.76817682345034,14.02618495842891,19.87072089461796,10.346029278056),
    @(10.96327628730096,6.24544514893531,5.486145951092141,23.76626250999215,27.77504559695381,14.02964336669727,19.87745485522089,10.33096481443218),
    @(10.89652891430259,6.200919706699294,5.473073324382853,23.7777256294434,27.79776897537491,14.04099374063607,19.89953796075639,10.28170252144081),
    @(10.85545132688382,6.173470219653917,5.465043200716157,23.78490888918247,27.81193538813881,14.04800130003901,19.91315838560757,10.2514261888433),
    @(10.61714530371809,6.013473969551909,5.418683886241793,23.82854688543188,27.89689837182621,14.08900339833818,19.9926500253383,10.07641005360234),
    @(10.46841500365434,5.91292781914113,5.389946526905828,23.85751816770782,27.95237738882149,14.11491014924452,20.04269863770182,9.967742260305659),
    @(10.381960500854,5.854221265571462,5.373311347845757,23.87498782903343,27.98551437546476,14.13008771786996,20.07195660197133,9.904782038746854),
    @(10.35253477144541,5.834194415366404,5.36766101780481,23.88104095692297,27.99694406875647,14.13527409074084,20.08194377923762,9.883388270355598),
    @(10.48434220589832,5.923721804704862,5.393016730640592,23.85435063172875,27.94634441744571,14.11212369180659,20.03732209274175,9.979357949400249),
    @(10.91617454978844,6.214034683527609,5.476917831902783,23.77432490732282,27.79104291389719,14.03764837127458,19.8930320505256,10.29619336215919),
    @(11.18982999813543,6.395888720783273,5.530741507348801,23.72925683478011,27.70061222178771,13.99144143838389,19.80293169240289,10.49876389067634),
    @(11.04458060368149,6.299554772671827,5.502110322206915,23.75264611503846,27.74785730396277,14.01587706053659,19.85063527454496,10.39107942381004),
    @(10.47714446335575,5.918844649932391,5.391629044913189,23.85578014078175,27.9490680586218,14.11338256635698,20.0397513200183,9.974108004184238),
    @(9.830713058971362,5.474818262937545,5.268438144341965,23.99780994608748,28.21292491252607,14.22913520100286,20.26176212465547,9.507160704929493)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $row).Value = $rowValues[$j]
    }
}
